$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1873
$ws1.Range("F7").Value = 3733
$ws1.Range("F13").Value = 654
$ws1.Range("F14").Value = 154
$ws1.Range("F15").Value = 870
$ws1.Range("F17").Value = 221
$ws1.Range("F22").Value = 3146
$ws1.Range("F23").Value = 5520
$ws1.Range("F28").Value = 3182
$ws1.Range("F30").Value = 2378
$ws1.Range("F32").Value = 507
$ws1.Range("F34").Value = 170
$ws1.Range("F37").Value = 87
$ws1.Range("F38").Value = 489
$ws1.Range("F39").Value = 858
$ws1.Range("F41").Value = 25
$ws1.Range("F44").Value = 526

# Sheet "全部类型" (fourth sheet) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1873
$ws4.Range("F7").Value = 3733
$ws4.Range("F14").Value = 654
$ws4.Range("F15").Value = 154
$ws4.Range("F16").Value = 870
$ws4.Range("F18").Value = 221
$ws4.Range("F22").Value = 80
$ws4.Range("F23").Value = 3146
$ws4.Range("F24").Value = 5520
$ws4.Range("F29").Value = 3182
$ws4.Range("F31").Value = 2378
$ws4.Range("F33").Value = 507
$ws4.Range("F35").Value = 170
$ws4.Range("F38").Value = 87
$ws4.Range("F39").Value = 489
$ws4.Range("F40").Value = 858
$ws4.Range("F42").Value = 25
$ws4.Range("F45").Value = 526
